# Applies the "Add files via upload" edit to the pbl_met description deck:
#  1. Inserts a new slide ("Changes since legacy pbl_met") at position 5.
#  2. Expands the "Test code" slide (now pushed to position 6) with extra bullets.
#  3. Retitles the "Why fortran?" slide (now position 7) to "But: Why fortran?".
#  4. Leaves "Fortran specific advantages" (now position 8) untouched.
#  5. Edits the "Code readability" slide (now position 9): fixes the
#     "philosophy" sentence and appends a final bullet.

function Set-Italic($paragraph, $start, $len) {
    $chars = $paragraph.Characters($start, $len)
    $chars.Font.Italic = $true
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. New slide at position 5: "Changes since legacy pbl_met"
# ---------------------------------------------------------------------
$refLayout = $p.Slides.Item(5).CustomLayout
$newSlide = $p.Slides.AddSlide(5, $refLayout)

$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Changes since legacy pbl_met"

$newBody = $newSlide.Shapes.Item(2).TextFrame.TextRange
$newBody.Text = "Focus expanded, from 1D met processors to general micro-meteorological data processing.`r" + `
"As a result, the code base enlarged significantly:`r" + `
"Statistics!`r" + `
"Eddy-covariance (wind, temperature, water (optional) and carbon dioxide (optional)).`r" + `
"Support for sensors like ultrasonic anemometers (various data formats), SODAR/RASS and narrow-beam disdrometric radars.`r" + `
"Test suites now part of code base.`r" + `
"Example applications provided.`r" + `
"In general, extensive refactoring of legacy code."

$newBody.Paragraphs(3).IndentLevel = 2
$newBody.Paragraphs(4).IndentLevel = 2
$newBody.Paragraphs(5).IndentLevel = 2

# ---------------------------------------------------------------------
# 2. "Test code" slide, now at position 6: add four new bullets.
# ---------------------------------------------------------------------
$testSlide = $p.Slides.Item(6)
$testBody = $testSlide.Shapes.Item(2).TextFrame.TextRange
$testBody.Text = "Quite a part of pbl_met is " + [char]8220 + "test code" + [char]8221 + ", aimed at harnessing procedure functionalities, and " + [char]8220 + "proving" + [char]8221 + " they are correct.`r" + `
"For important modules, test code size equals or exceeds module size.`r" + `
"This is unlike the legacy PBL_MET, for which test code was not released.`r" + `
"Test code, collected under directory " + [char]8220 + "/test" + [char]8221 + ", also provides examples on using individual routines.`r" + `
"Test code also useful to see various modern Fortran constructs at work.`r" + `
"Named optional parameters in argument list.`r" + `
"Object-orientation."

Set-Italic $testBody.Paragraphs(1) 17 7

$testBody.Paragraphs(2).IndentLevel = 2
$testBody.Paragraphs(6).IndentLevel = 2
$testBody.Paragraphs(7).IndentLevel = 2

# ---------------------------------------------------------------------
# 3. "Why fortran?" slide, now at position 7: retitle only.
# ---------------------------------------------------------------------
$whySlide = $p.Slides.Item(7)
$whySlide.Shapes.Item(1).TextFrame.TextRange.Text = "But: Why fortran?"

# ---------------------------------------------------------------------
# 4. "Fortran specific advantages" slide, now at position 8: untouched.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 5. "Code readability" slide, now at position 9: edit + append bullet.
# ---------------------------------------------------------------------
$readSlide = $p.Slides.Item(9)
$readBody = $readSlide.Shapes.Item(2).TextFrame.TextRange

$readBody.Paragraphs(2).Text = "Our philosophy departs somewhat from open-source as traditionally intended in meteorology and atmospheric physics: we firmly believe there is no logical reason people, in front of a " + [char]8220 + "public" + [char]8221 + " project, should " + [char]8220 + "sink or swim" + [char]8221 + "."

$readBody.InsertAfter("`rAnd, high value placed for code whose meaning is transparent, even if this may decrease efficiency to some extent.")
